# Add new Courts-related test-data mapping rows to the
# TestDataMappingSheet_SD sheet, continuing the pattern of the rows
# directly above (ScreenName / TestDataFileName / TestDataSheetName /
# StartIndexofIteration columns tied to cares\Courts.xlsx).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestDataMappingSheet_SD")
$ws.Activate()

# Rows 134-135 already exist (blank); insert 6 more formatted rows below
# them (136-141) so every new row picks up the same look as the rows
# already in the table.
$ws.Range("A136:A141").EntireRow.Insert()

# Row 141 is left as a single trailing blank/formatted cell, like the
# blank row that trails the pasted block in the source file.
$ws.Range("E136:E140").Clear()
$ws.Range("A141").Clear()
$ws.Range("C141:E141").Clear()

$testDataFile = "cares\Courts.xlsx"

$newRows = @(
    @(134, "CourtOfficerNotes"),
    @(135, "CourtReports"),
    @(136, "CourtMinuteOrders"),
    @(137, "DocumentDistributions"),
    @(138, "RelatedFolios"),
    @(139, "CourtsApprovalAndAuditHistory"),
    @(140, "CourtsAuditHistory")
)

foreach ($entry in $newRows) {
    $r = $entry[0]
    $screenName = $entry[1]

    $ws.Cells.Item($r, 1).Value = $screenName
    $ws.Cells.Item($r, 2).Value = $testDataFile
    $ws.Cells.Item($r, 3).Value = $screenName
    $ws.Cells.Item($r, 4).Value = 1
}

# Reflect the author's final scroll position/selection on the sheet.
$excel.ActiveWindow.ScrollRow = 120
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("A130").Select()
